# Automatische test-sync: 2025-08-28 21:14:50
#
# Adds the new "CE-certificaat aanvraag" mail-log entry (Logs row 25) and
# its corresponding Dashboard aggregate row (row 6, "Kwaliteit /
# Certificaten" = 1), then extends the conditional formatting ranges on
# the Logs sheet and the bar chart's category/value series references on
# the Dashboard sheet so they include the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 25
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A25").Value = "CE-certificaat aanvraag"
$logs.Range("B25").Value = "mailmind.test@zohomail.eu"
$logs.Range("D25").Value = "Kwaliteit / Certificaten"
$logs.Range("F25").Value = "2025-08-28 21:14:27"
$logs.Range("G25").Value = "Nee"
$logs.Range("H25").Value = "Ja"
$logs.Range("I25").Value = "Nee"
$logs.Range("J25").Value = "Nee"

# Extend the conditional formatting sqref on the Logs sheet from row 24
# to row 25 for every column that carries a rule (D, G, H, I, J).
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range($col + "2:" + $col + "24")
    $newRange = $logs.Range($col + "2:" + $col + "25")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet: append aggregate row 6
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A6").Value = "Kwaliteit / Certificaten"
$dash.Range("B6").Value = 1

# ---------------------------------------------------------------------
# 3. Chart on the Dashboard sheet: extend the series cat/val references
#    from $A$2:$A$5 / $B$2:$B$5 to $A$2:$A$6 / $B$2:$B$6
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects(1).Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$6,Dashboard!`$B`$2:`$B`$6,1)"
